$wb = $excel.ActiveWorkbook

# Sheet ALC, row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 949.2308
$ws.Range("I118").Value = 356.66666
$ws.Range("J118").Value = 1457.1428
$ws.Range("K118").Value = 1069.99998
$ws.Range("L118").Value = 4371.428400000001
$ws.Range("M118").Value = 587.0000199999999
$ws.Range("N118").Value = -7685.428400000001

# Sheet ALC, row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2091.842
$ws.Range("I127").Value = 605.7143
$ws.Range("J127").Value = 2958.75
$ws.Range("K127").Value = 1817.1429
$ws.Range("L127").Value = 8876.25
$ws.Range("M127").Value = 3142.8571
$ws.Range("N127").Value = -18796.25

# Sheet ALC, row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1860.1177
$ws.Range("I131").Value = 1100.5385
$ws.Range("K131").Value = 3301.6155
$ws.Range("M131").Value = 1738.3845

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 852.5925999999999
$ws.Range("I137").Value = 747.25
$ws.Range("J137").Value = 936.86664
$ws.Range("K137").Value = 2241.75
$ws.Range("L137").Value = 2810.59992
$ws.Range("M137").Value = 308.25
$ws.Range("N137").Value = -7910.59992

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1534.2
$ws.Range("I138").Value = 689.371
$ws.Range("J138").Value = 2912.6052
$ws.Range("K138").Value = 2068.113
$ws.Range("L138").Value = 8737.8156
$ws.Range("M138").Value = 3071.887
$ws.Range("N138").Value = -19017.8156

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2012.1086
$ws.Range("I141").Value = 746.64514
$ws.Range("J141").Value = 4627.4
$ws.Range("K141").Value = 2239.93542
$ws.Range("L141").Value = 13882.2
$ws.Range("M141").Value = 2940.06458
$ws.Range("N141").Value = -24242.2

# Sheet ARM, row 25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 5508
$ws.Range("I25").Value = 1016
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 1016
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = -614
$ws.Range("N25").Value = -10804

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19866.918
$ws.Range("I32").Value = 23091.61
$ws.Range("J32").Value = 10039.286
$ws.Range("K32").Value = 23091.61
$ws.Range("L32").Value = 10039.286
$ws.Range("M32").Value = -22804.61
$ws.Range("N32").Value = -10613.286

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 881.7059
$ws.Range("I61").Value = 694.60974
$ws.Range("J61").Value = 1648.8
$ws.Range("K61").Value = 694.60974
$ws.Range("L61").Value = 1648.8
$ws.Range("M61").Value = -482.60974
$ws.Range("N61").Value = -2072.8

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 821.6875
$ws.Range("I110").Value = 749.7692
$ws.Range("K110").Value = 749.7692
$ws.Range("M110").Value = 1295.2308

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1022.75
$ws.Range("I132").Value = 844.3137
$ws.Range("J132").Value = 2842.8
$ws.Range("K132").Value = 2532.9411
$ws.Range("L132").Value = 8528.400000000001
$ws.Range("M132").Value = -2.941100000000006
$ws.Range("N132").Value = -13588.4

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 881.7059
$ws.Range("I136").Value = 694.60974
$ws.Range("J136").Value = 1648.8
$ws.Range("K136").Value = 2083.82922
$ws.Range("L136").Value = 4946.4
$ws.Range("M136").Value = 466.1707799999999
$ws.Range("N136").Value = -10046.4

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1184.8214
$ws.Range("I94").Value = 1189.3182
$ws.Range("J94").Value = 1168.3334
$ws.Range("K94").Value = 1189.3182
$ws.Range("L94").Value = 1168.3334
$ws.Range("M94").Value = -738.3181999999999
$ws.Range("N94").Value = -2070.3334

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1829.9474
$ws.Range("I99").Value = 1847.8572
$ws.Range("J99").Value = 1779.8
$ws.Range("K99").Value = 1847.8572
$ws.Range("L99").Value = 1779.8
$ws.Range("M99").Value = -349.8571999999999
$ws.Range("N99").Value = -4775.8

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14963.054
$ws.Range("I134").Value = 1161.9844
$ws.Range("J134").Value = 103289.9
$ws.Range("K134").Value = 3485.9532
$ws.Range("L134").Value = 309869.7
$ws.Range("M134").Value = -950.9531999999999
$ws.Range("N134").Value = -314939.7

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2423.2856
$ws.Range("I31").Value = 2551.5715
$ws.Range("J31").Value = 2209.476
$ws.Range("K31").Value = 2551.5715
$ws.Range("L31").Value = 2209.476
$ws.Range("M31").Value = -2256.5715
$ws.Range("N31").Value = -2799.476

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2423.2856
$ws.Range("I34").Value = 2551.5715
$ws.Range("J34").Value = 2209.476
$ws.Range("K34").Value = 2551.5715
$ws.Range("L34").Value = 2209.476
$ws.Range("M34").Value = -2349.5715
$ws.Range("N34").Value = -2613.476

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 961.1667
$ws.Range("I5").Value = 1037.1111
$ws.Range("J5").Value = 733.3333
$ws.Range("K5").Value = 3111.3333
$ws.Range("L5").Value = 2199.9999
$ws.Range("M5").Value = -2999.3333
$ws.Range("N5").Value = -2423.9999

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1297369.4
$ws.Range("I107").Value = 1766.3334
$ws.Range("J107").Value = 2592972.2
$ws.Range("K107").Value = 5299.0002
$ws.Range("L107").Value = 7778916.600000001
$ws.Range("M107").Value = -3379.0002
$ws.Range("N107").Value = -7782756.600000001

# Sheet CUL, row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2707.8
$ws.Range("I109").Value = 908.5
$ws.Range("J109").Value = 2984.6155
$ws.Range("K109").Value = 2725.5
$ws.Range("L109").Value = 8953.8465
$ws.Range("M109").Value = -1685.5
$ws.Range("N109").Value = -11033.8465

# Sheet CUL, row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1022.63635
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1022.63635
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").Value = 3067.90905
$ws.Range("N121").Value = -5687.90905

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 263562.34
$ws.Range("I122").Value = 178.83333
$ws.Range("J122").Value = 500607.5
$ws.Range("K122").Value = 1609.49997
$ws.Range("L122").Value = 4505467.5
$ws.Range("M122").Value = 840.5000300000002
$ws.Range("N122").Value = -4510367.5

# Sheet CUL, row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 5450
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 5450
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").Value = 16350
$ws.Range("N125").Value = -26190

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 961.1667
$ws.Range("I135").Value = 1037.1111
$ws.Range("J135").Value = 733.3333
$ws.Range("K135").Value = 9333.999900000001
$ws.Range("L135").Value = 6599.9997
$ws.Range("M135").Value = -6798.999900000001
$ws.Range("N135").Value = -11669.9997

# Sheet GSM, row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 7021002
$ws.Range("J24").Value = 70007
$ws.Range("L24").Value = 70007
$ws.Range("N24").Value = -70353

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2014.8857
$ws.Range("I132").Value = 1823.258
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 5469.774
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -2939.774
$ws.Range("N132").Value = -15560

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2633.3333
$ws.Range("I100").Value = 3366.6667
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 3366.6667
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -2825.6667
$ws.Range("N100").Value = -2982

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2228.1843
$ws.Range("I136").Value = 1240.7222
$ws.Range("J136").Value = 20002.5
$ws.Range("K136").Value = 3722.1666
$ws.Range("L136").Value = 60007.5
$ws.Range("M136").Value = -1172.1666
$ws.Range("N136").Value = -65107.5

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 655.40625
$ws.Range("I132").Value = 528.65454
$ws.Range("J132").Value = 1430
$ws.Range("K132").Value = 1585.96362
$ws.Range("L132").Value = 4290
$ws.Range("M132").Value = 944.03638
$ws.Range("N132").Value = -9350
